$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the target cell and its new value, taken from the
# "after" side of the authoritative diff. Columns D (Price) and E
# (Volume(1h)) are stored as text in the workbook (e.g. "65.369.64",
# "  -3.73%  "), so NumberFormat is forced to "@" (Text) before the
# assignment to stop Excel from re-interpreting numeric-looking
# strings (and silently dropping trailing zeros, e.g. "1.00" -> 1).
$updates = @(
    @{ Cell = 'D2'; Value = '65.369.64' }
    @{ Cell = 'E2'; Value = '  -3.73%  ' }
    @{ Cell = 'D3'; Value = '3.484.83' }
    @{ Cell = 'E3'; Value = '  -1.11%  ' }
    @{ Cell = 'E4'; Value = '  +0.08%  ' }
    @{ Cell = 'D5'; Value = '552.96' }
    @{ Cell = 'E5'; Value = '  -0.76%  ' }
    @{ Cell = 'D6'; Value = '179.07' }
    @{ Cell = 'E6'; Value = '  -6.71%  ' }
    @{ Cell = 'E7'; Value = '  +4.44%  ' }
    @{ Cell = 'E8'; Value = '  +0.02%  ' }
    @{ Cell = 'D9'; Value = '0.632' }
    @{ Cell = 'E9'; Value = '  -1.30%  ' }
    @{ Cell = 'D10'; Value = '0.155' }
    @{ Cell = 'E10'; Value = '  +2.97%  ' }
    @{ Cell = 'D11'; Value = '53.89' }
    @{ Cell = 'E11'; Value = '  -5.88%  ' }
    @{ Cell = 'E12'; Value = '  -2.25%  ' }
    @{ Cell = 'D13'; Value = '9.16' }
    @{ Cell = 'E13'; Value = '  -3.57%  ' }
    @{ Cell = 'D14'; Value = '4.048.22' }
    @{ Cell = 'E14'; Value = '  -0.96%  ' }
    @{ Cell = 'D15'; Value = '3.491.20' }
    @{ Cell = 'E15'; Value = '  -0.95%  ' }
    @{ Cell = 'E16'; Value = '  -0.05%  ' }
    @{ Cell = 'D17'; Value = '18.37' }
    @{ Cell = 'D18'; Value = '12.18' }
    @{ Cell = 'E18'; Value = '  +2.16%  ' }
    @{ Cell = 'D19'; Value = '65.462.99' }
    @{ Cell = 'E19'; Value = '  -3.55%  ' }
    @{ Cell = 'D20'; Value = '0.994' }
    @{ Cell = 'E20'; Value = '  -1.61%  ' }
    @{ Cell = 'D21'; Value = '413.02' }
    @{ Cell = 'E21'; Value = '  +1.23%  ' }
    @{ Cell = 'D22'; Value = '4.04' }
    @{ Cell = 'E22'; Value = '  +1.91%  ' }
    @{ Cell = 'E23'; Value = '  +0.87%  ' }
    @{ Cell = 'E24'; Value = '  -3.38%  ' }
    @{ Cell = 'D25'; Value = '12.73' }
    @{ Cell = 'E25'; Value = '  +6.72%  ' }
    @{ Cell = 'D26'; Value = '10.78' }
    @{ Cell = 'E26'; Value = '  -8.03%  ' }
    @{ Cell = 'D27'; Value = '2.85' }
    @{ Cell = 'E27'; Value = '  -2.17%  ' }
    @{ Cell = 'B28'; Value = 'Filecoin' }
    @{ Cell = 'C28'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' }
    @{ Cell = 'D28'; Value = '9.00' }
    @{ Cell = 'E28'; Value = '  +4.19%  ' }
    @{ Cell = 'B29'; Value = 'EthereumClassic' }
    @{ Cell = 'C29'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' }
    @{ Cell = 'D29'; Value = '30.26' }
    @{ Cell = 'E29'; Value = '  -1.10%  ' }
    @{ Cell = 'B30'; Value = 'Bittensor' }
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao' }
    @{ Cell = 'D30'; Value = '612.93' }
    @{ Cell = 'E30'; Value = '  -10.06%  ' }
    @{ Cell = 'B31'; Value = 'NEARProtocol' }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' }
    @{ Cell = 'D31'; Value = '6.45' }
    @{ Cell = 'E31'; Value = '  -6.46%  ' }
    @{ Cell = 'B32'; Value = 'Cosmos' }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' }
    @{ Cell = 'D32'; Value = '11.65' }
    @{ Cell = 'E32'; Value = '  -0.80%  ' }
    @{ Cell = 'B33'; Value = 'Hedera' }
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' }
    @{ Cell = 'D33'; Value = '0.110' }
    @{ Cell = 'E33'; Value = '  -1.42%  ' }
    @{ Cell = 'B34'; Value = 'OKB' }
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb' }
    @{ Cell = 'D34'; Value = '59.49' }
    @{ Cell = 'E34'; Value = '  -1.64%  ' }
    @{ Cell = 'B35'; Value = 'Kaspa' }
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas' }
    @{ Cell = 'D35'; Value = '0.148' }
    @{ Cell = 'E35'; Value = '  +11.01%  ' }
    @{ Cell = 'B36'; Value = 'Dai' }
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' }
    @{ Cell = 'D36'; Value = '1.00' }
    @{ Cell = 'E36'; Value = '  +0.10%  ' }
    @{ Cell = 'B37'; Value = 'PEPE' }
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe' }
    @{ Cell = 'D37'; Value = '0.0₃0789' }
    @{ Cell = 'E37'; Value = '  -7.77%  ' }
    @{ Cell = 'D38'; Value = '37.06' }
    @{ Cell = 'E38'; Value = '  -5.48%  ' }
    @{ Cell = 'B39'; Value = 'Maker' }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr' }
    @{ Cell = 'D39'; Value = '3.366.53' }
    @{ Cell = 'E39'; Value = '  +10.19%  ' }
    @{ Cell = 'B40'; Value = 'TheGraph' }
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt' }
    @{ Cell = 'D40'; Value = '0.380' }
    @{ Cell = 'E40'; Value = '  -6.18%  ' }
    @{ Cell = 'B41'; Value = 'FirstDigitalUSD' }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd' }
    @{ Cell = 'D41'; Value = '0.999' }
    @{ Cell = 'E41'; Value = '  -0.06%  ' }
    @{ Cell = 'B42'; Value = 'Stacks' }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx' }
    @{ Cell = 'D42'; Value = '3.25' }
    @{ Cell = 'E42'; Value = '  -4.89%  ' }
    @{ Cell = 'B43'; Value = 'ThetaToken' }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta' }
    @{ Cell = 'D43'; Value = '2.84' }
    @{ Cell = 'E43'; Value = '  -5.95%  ' }
    @{ Cell = 'B44'; Value = 'WEMIXToken' }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' }
    @{ Cell = 'D44'; Value = '2.77' }
    @{ Cell = 'E44'; Value = '  +0.81%  ' }
    @{ Cell = 'E45'; Value = '  -9.62%  ' }
    @{ Cell = 'B46'; Value = 'VeChain' }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Cell = 'D46'; Value = '0.0414' }
    @{ Cell = 'E46'; Value = '  -2.20%  ' }
    @{ Cell = 'B47'; Value = 'ApeXProtocol' }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex' }
    @{ Cell = 'D47'; Value = '3.23' }
    @{ Cell = 'E47'; Value = '  +0.66%  ' }
    @{ Cell = 'B48'; Value = 'Stellar' }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm' }
    @{ Cell = 'D48'; Value = '0.133' }
    @{ Cell = 'E48'; Value = '  +1.44%  ' }
    @{ Cell = 'B49'; Value = 'THORChain' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune' }
    @{ Cell = 'D49'; Value = '8.42' }
    @{ Cell = 'E49'; Value = '  -10.04%  ' }
    @{ Cell = 'D50'; Value = '137.36' }
    @{ Cell = 'E50'; Value = '  -1.10%  ' }
    @{ Cell = 'B51'; Value = 'LidoDAOToken' }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' }
    @{ Cell = 'D51'; Value = '2.86' }
    @{ Cell = 'E51'; Value = '  +9.65%  ' }
)

foreach ($update in $updates) {
    $range = $ws.Range($update.Cell)
    if ($update.Cell[0] -eq 'D' -or $update.Cell[0] -eq 'E') {
        $range.NumberFormat = '@'
    }
    $range.Value = $update.Value
}
